$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear values in C2:C13 (keep style/number format)
$ws.Range("C2:C13").ClearContents()

# Update selection to match the diff (C2 active cell, C2:C13 selected)
$ws.Range("C2:C13").Select()
